$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.472.38"
$ws.Range("E2").Value = "  -1.62%  "

$ws.Range("D3").Value = "3.081.32"
$ws.Range("E3").Value = "  -2.70%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.78"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.50"
$ws.Range("E6").Value = "  +3.23%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("E8").Value = "  +5.68%  "

$ws.Range("D9").Value = "3.079.92"
$ws.Range("E9").Value = "  -2.49%  "

$ws.Range("E10").Value = "  -4.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  -0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("E12").Value = "  -0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000242"
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.35"
$ws.Range("E14").Value = "  -0.58%  "

$ws.Range("E15").Value = "  -1.92%  "

$ws.Range("D16").Value = "3.591.71"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.22"
$ws.Range("E17").Value = "  -0.59%  "

$ws.Range("D18").Value = "63.504.68"
$ws.Range("E18").Value = "  -1.29%  "

$ws.Range("D19").Value = "3.082.88"
$ws.Range("E19").Value = "  -2.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.06"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"
$ws.Range("E21").Value = "  +1.03%  "

$ws.Range("E22").Value = "  -2.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.11"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("E25").Value = "  -0.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.81"
$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  -0.81%  "

$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").Value = "  -3.11%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.25"
$ws.Range("E32").Value = "  -2.99%  "

$ws.Range("E33").Value = "  +1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.39"
$ws.Range("E34").Value = "  -1.16%  "

$ws.Range("D35").Value = "0.0₃0855"
$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("E36").Value = "  -1.81%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.14"
$ws.Range("E37").Value = "  -2.18%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +2.72%  "

$ws.Range("E39").Value = "  -5.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.39"
$ws.Range("E40").Value = "  +1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.48"
$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "443.03"
$ws.Range("E42").Value = "  -6.35%  "

$ws.Range("E43").Value = "  -4.83%  "

$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "39.88"
$ws.Range("E45").Value = "  -2.91%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.815.65"
$ws.Range("E46").Value = "  -4.13%  "

$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.110"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.54"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E50").Value = "  +1.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.01"
$ws.Range("E51").Value = "  +1.69%  "
